# Applies crypto price/volume updates to sheet1, matching the commit
# "Updated symbol list on Tue Jan 24 03:26:22 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "319.26"
Set-TextValue $ws.Range("E2") "4.44%"

# Row 3
Set-TextValue $ws.Range("D3") "35.99"
Set-TextValue $ws.Range("E3") "-0.63%"

# Row 4
Set-TextValue $ws.Range("D4") "5.123"
Set-TextValue $ws.Range("E4") "1.78%"

# Row 5
Set-TextValue $ws.Range("D5") "0.08060"
Set-TextValue $ws.Range("E5") "2.98%"

# Row 6
Set-TextValue $ws.Range("D6") "2.163"
Set-TextValue $ws.Range("E6") "2.04%"

# Row 7
Set-TextValue $ws.Range("D7") "8.059"
Set-TextValue $ws.Range("E7") "1.72%"

# Row 8
Set-TextValue $ws.Range("D8") "4.130"
Set-TextValue $ws.Range("E8") "1.65%"

# Row 9
Set-TextValue $ws.Range("D9") "0.9290"
Set-TextValue $ws.Range("E9") "1.16%"

# Row 10
Set-TextValue $ws.Range("D10") "0.1009"
Set-TextValue $ws.Range("E10") "5.42%"

# Row 11
Set-TextValue $ws.Range("E11") "0.55%"

# Row 12
Set-TextValue $ws.Range("D12") "0.09196"
Set-TextValue $ws.Range("E12") "5.71%"

# Row 13
Set-TextValue $ws.Range("D13") "0.03598"
Set-TextValue $ws.Range("E13") "2.61%"

# Row 14
Set-TextValue $ws.Range("D14") "0.09931"

# Row 15
Set-TextValue $ws.Range("D15") "0.001432"
Set-TextValue $ws.Range("E15") "-0.50%"

# Row 16
Set-TextValue $ws.Range("D16") "0.005661"
Set-TextValue $ws.Range("E16") "-0.95%"

# Row 17
Set-TextValue $ws.Range("D17") "3.462"
Set-TextValue $ws.Range("E17") "0.09%"

# Row 18
Set-TextValue $ws.Range("E18") "18.43%"

# Row 19
Set-TextValue $ws.Range("E19") "-1.25%"

# Row 20
Set-TextValue $ws.Range("D20") "0.1330"
Set-TextValue $ws.Range("E20") "-1.14%"

# Row 21
Set-TextValue $ws.Range("D21") "5.059"
Set-TextValue $ws.Range("E21") "6.00%"

# Row 22
Set-TextValue $ws.Range("D22") "0.2204"
Set-TextValue $ws.Range("E22") "-3.80%"

# Row 23
Set-TextValue $ws.Range("D23") "0.04614"
Set-TextValue $ws.Range("E23") "0.14%"

# Row 24
Set-TextValue $ws.Range("D24") "0.001244"
Set-TextValue $ws.Range("E24") "1.10%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004744"
Set-TextValue $ws.Range("E25") "-6.99%"

# Row 26
Set-TextValue $ws.Range("D26") "0.0001302"
Set-TextValue $ws.Range("E26") "-7.11%"

# Row 27
Set-TextValue $ws.Range("D27") "0.0004504"
Set-TextValue $ws.Range("E27") "65.22%"

# Row 39
Set-TextValue $ws.Range("D39") "0.01960"
Set-TextValue $ws.Range("E39") "7.19%"

# Row 40
Set-TextValue $ws.Range("D40") "0.04979"
Set-TextValue $ws.Range("E40") "4.43%"

# Row 41
Set-TextValue $ws.Range("D41") "0.007812"
Set-TextValue $ws.Range("E41") "4.23%"

# Row 42
Set-TextValue $ws.Range("D42") "0.1400"
Set-TextValue $ws.Range("E42") "-0.04%"

# Row 43
Set-TextValue $ws.Range("D43") "0.007630"
Set-TextValue $ws.Range("E43") "-1.28%"

# Row 44
Set-TextValue $ws.Range("D44") "0.002069"
Set-TextValue $ws.Range("E44") "-7.30%"

# Row 45
Set-TextValue $ws.Range("D45") "0.01170"
Set-TextValue $ws.Range("E45") "12.05%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006312"
Set-TextValue $ws.Range("E46") "1.32%"

# Row 47
Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "0.03%"

# Row 48
Set-TextValue $ws.Range("E48") "2.52%"

# Row 49
Set-TextValue $ws.Range("D49") "0.001902"
Set-TextValue $ws.Range("E49") "-4.96%"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("E50") "0.03%"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002002"
Set-TextValue $ws.Range("E51") "0.03%"

